$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replace existing account holder's data (LEVI -> ANA).
# Force the account-number column to text so leading zeros survive.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "004210959"
$ws.Range("B5").Value = "ANA"
$ws.Range("C5").Value = 285000

# Insert a new row after row 6 (the MARCUS / 000834301 row) for DOUGLAS,
# shifting the remaining rows down.
$ws.Rows.Item(7).Insert()
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "004384167"
$ws.Range("B7").Value = "DOUGLAS"
$ws.Range("C7").Value = 87159.6
